{"js": "// Remove now-redundant empty spacer paragraphs (\"NodeEnd\" / \"HeadEnd\")\n// that were leaving blank separator lines around headers with link\n// resources. Every \"Node End\" spacer paragraph is dropped, and a\n// \"Head End\" spacer paragraph is dropped whenever it sits directly in\n// front of a \"Body Text\" paragraph (i.e. whenever the heading actually\n// has body content following it, the extra blank line is unnecessary).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"style\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst indexesToRemove = [];\n\nfor (let i = 0; i < items.length; i++) {\n  const style = items[i].style;\n  if (style === \"Node End\") {\n    indexesToRemove.push(i);\n  } else if (\n    style === \"Head End\" &&\n    i + 1 < items.length &&\n    items[i + 1].style === \"Body Text\"\n  ) {\n    indexesToRemove.push(i);\n  }\n}\n\n// Delete from the highest index down so earlier indexes stay valid.\nindexesToRemove.sort((a, b) => b - a);\nfor (const idx of indexesToRemove) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove now-redundant empty spacer paragraphs (\"NodeEnd\" / \"HeadEnd\")\n# that were leaving blank separator lines around headers with link\n# resources. Every \"Node End\" spacer paragraph is dropped, and a\n# \"Head End\" spacer paragraph is dropped whenever it sits directly in\n# front of a \"Body Text\" paragraph (i.e. whenever the heading actually\n# has body content following it, the extra blank line is unnecessary).\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# Snapshot every paragraph's style name up front so later deletions\n# don't perturb the indices we still need to inspect.\n$styles = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $styles += $d.Paragraphs.Item($i).Style.NameLocal\n}\n\n$toRemove = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $style = $styles[$i - 1]\n    if ($style -eq \"Node End\") {\n        $toRemove += $i\n    } elseif ($style -eq \"Head End\" -and $i -lt $count -and $styles[$i] -eq \"Body Text\") {\n        $toRemove += $i\n    }\n}\n\n# Delete from the highest index down so earlier indexes stay valid.\n$sorted = $toRemove | Sort-Object -Descending\nforeach ($idx in $sorted) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
